$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "Sun Jan 14 16:33:24 EST 2024"
$ws.Range("B3").Value = "Sun Jan 14 16:33:36 EST 2024"
$ws.Range("B5").Value = "Sun Jan 14 16:33:47 EST 2024"
